$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K6").Value = "69f45ebf-bc69-444e-9e59-72e123d1afa4"
$ws.Range("K7").Value = "dbaf647a-d05a-455c-9ef3-94131b1a2f22"
$ws.Range("K8").Value = "fc57ae5a-8d70-4b43-bb4c-0a513397997c"
$ws.Range("K9").Value = "7d13ba6b-d02a-4a8c-b65a-83ca38438f05"
$ws.Range("K10").Value = "eed40001-1933-455a-8e39-7e3892ed2ce5"
$ws.Range("K11").Value = "9b3b923a-5646-427f-b0a3-18cf52bb974d"
$ws.Range("K12").Value = "f25dcc11-b4bc-4ae5-90da-17a8cb58c6eb"
$ws.Range("K13").Value = "6ae02d52-1294-4e4c-9a4e-d93c28e27e9b"
$ws.Range("K14").Value = "850d23cc-ae5b-4e96-ad78-ff0d692ed313"
$ws.Range("K15").Value = "e1165481-4678-48dd-8072-2d5a8d8994a6"
$ws.Range("K16").Value = "4dc156fa-eada-4fa7-a2f3-5aede98788eb"
$ws.Range("K17").Value = "d86cd5a3-fa2d-4804-9fcb-81e349e4b166"
$ws.Range("G18").Value = "b41ccd27-9a4f-5cc8-9c5d-b55242d90fb0"
$ws.Range("K18").Value = "7ef8bdeb-fd56-5eb9-a09b-ef15ce18dc49"
$ws.Range("G19").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G20").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G21").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G22").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G23").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G24").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G25").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G26").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G27").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G28").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G29").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G30").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G31").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G32").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G33").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G34").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G35").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G36").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G37").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G38").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("G39").Value = "591191c7-f693-5957-8734-ac87151ca981"
$ws.Range("K19").Value = "ed7b1af2-3c5c-4259-a997-5eaeca20f2d5"
$ws.Range("K20").Value = "815b48f7-b310-4ff6-aa02-6b8f7c049bb1"
$ws.Range("K21").Value = "53dd6491-04c6-45f1-abf7-55cb5180a256"
$ws.Range("K22").Value = "5846bdfe-fb03-4f05-87e7-641a2b943ff7"
$ws.Range("K23").Value = "54d572f5-f901-4688-8a27-c4f3db636f11"
$ws.Range("K24").Value = "9ad84ef6-184a-437f-8472-31563a623011"
$ws.Range("K25").Value = "494294b4-1d59-4221-9fa6-0ba460a71409"
$ws.Range("K26").Value = "5225b69e-1ecc-4362-bd59-fa2bcc642317"
$ws.Range("K27").Value = "a5ede658-f75d-4d73-9297-962e1f4b7d72"
$ws.Range("K28").Value = "5e59c68e-51d4-4a57-bbca-ff8ec3ebe05f"
$ws.Range("K29").Value = "90ea63e8-a640-4b01-b573-e2924fc5f0ad"
$ws.Range("K30").Value = "e8d424f3-5833-4603-819d-23cc9039c289"
$ws.Range("K31").Value = "76865ebd-8524-4549-b034-3cbe130beb18"
$ws.Range("K32").Value = "f14ca97a-4bf1-46d4-8ca2-2eadf9e16556"
$ws.Range("K33").Value = "07eedb66-f428-4be7-a42c-ae48c4c28c20"
$ws.Range("K34").Value = "0dd7f6a8-0e76-4008-a21e-4a0a5de0ebc3"
$ws.Range("K35").Value = "6a9f8a4c-e08b-447d-9f3d-8d2ddd34979f"
$ws.Range("K36").Value = "868b5676-29a9-48d8-ba3a-feed82e9caac"
$ws.Range("K37").Value = "59f3777c-bb08-442c-816b-9727ccfda685"
$ws.Range("K38").Value = "011c4358-b74a-42d7-b4ad-586c1217a7a3"
$ws.Range("K39").Value = "1a8c1aad-0a8f-4bf0-9306-592948d513dd"
$ws.Range("G40").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G41").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G42").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G43").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G44").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G45").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G46").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G47").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G48").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G49").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G50").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G51").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G52").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G53").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G54").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G55").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G56").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G57").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G58").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G59").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G60").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G61").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G62").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G63").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G64").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G65").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("G66").Value = "86b4a49e-7378-5159-9f41-b005208c31bc"
$ws.Range("K40").Value = "fdee0ae4-6ff7-442e-b183-03041be75e7e"
$ws.Range("K41").Value = "8f431b0e-0d9e-42ca-97ca-7b61d8824d20"
$ws.Range("K42").Value = "2e79a982-fcdb-48c1-8752-ceee1aa77256"
$ws.Range("K43").Value = "a9c9e4e4-7828-41bc-8fc2-d7bb69ccf364"
$ws.Range("K44").Value = "2e33b845-3409-4bbf-9726-ee86617aa123"
$ws.Range("K45").Value = "7f952548-16c3-43c2-a1e5-e24275633ca4"
$ws.Range("K46").Value = "5f0237e5-b3c7-4951-842b-29f2ac710378"
$ws.Range("K47").Value = "e858541e-93ec-4ab1-bb80-870030c8de10"
$ws.Range("K48").Value = "e3e5a9eb-4f9a-41d9-8384-94781be61ce8"
$ws.Range("K49").Value = "f9b91d43-617f-4bd6-908e-e25fbc4a09ca"
$ws.Range("K50").Value = "48bfdaba-4762-477f-a736-b6048f4d54c1"
$ws.Range("K51").Value = "489580e5-e745-46fb-a8e9-1066b955698f"
$ws.Range("K52").Value = "13d70913-2ce9-4720-bb74-8e3e92266312"
$ws.Range("K53").Value = "fdf7d612-4cfb-44c4-a1bc-05706f66d404"
$ws.Range("K54").Value = "dd03acfd-5c0e-4411-be5c-92b4d78c5ec4"
$ws.Range("K55").Value = "1f629815-2287-49ab-932b-077a154b1f97"
$ws.Range("K56").Value = "4b5ede2d-7972-4025-b264-68ca35520a9d"
$ws.Range("K57").Value = "c9bba00c-0c5a-4be2-b6ae-2a2a94c7578b"
$ws.Range("K58").Value = "864968c5-c0f1-4a0c-9671-c0bb49a84693"
$ws.Range("K59").Value = "a5a593ff-757e-48db-a620-1951939f9e70"
$ws.Range("K60").Value = "3e5aeeef-7245-4a01-863b-b91e1341fc49"
$ws.Range("K61").Value = "04afe5bb-edee-4c11-a2e0-daf54d508d28"
$ws.Range("K62").Value = "da8a2660-535d-4430-b696-95f668ade476"
$ws.Range("K63").Value = "aaacb190-7de5-4069-b422-6a3184f9e631"
$ws.Range("K64").Value = "50d9b2ac-e7cd-4047-97f5-6c05d9bccd03"
$ws.Range("K65").Value = "afc2bc99-bc4f-450e-800c-ccfe303e27dd"
$ws.Range("K66").Value = "baacb865-a98a-4216-bfa0-8909c50afe23"
